$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StartSceneConfig")

# New row 13: mirror the formatting of row 12 (Account) then overwrite with
# the new "LoginCenter" entry. Copy propagates cell styles (and the H column
# entry, which is then cleared since row 13 has no H13 value).
$ws.Range("C12:H12").Copy($ws.Range("C13:H13"))

$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = "LoginCenter"
$ws.Range("G13").Value = "LoginCenter"
$ws.Range("H13").Clear()

# Update the view: scroll so row 4 is the top-left row, and move the
# selection to J11.
$ws.Activate()
$ws.Range("B4").Select()
$ws.Range("J11").Select()
